$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 149; $row++) {
    $ws.Cells.Item($row, 32).Value = "30/01/2026 17:47:34"
}
